# Daily attendance processing - 2025-12-31 20:58:40
# Swap the order of the recorder names in the "Recorded By" column (G)
# from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every session row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
